$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels in row 1 (processed D->A so the rebuilt shared-string
# table lands in the same order as the target: units_for_quantities_being_observed,
# type_of_quantity_observed, devicetype, deviceid).
$ws.Range("D1").Value = "units_for_quantities_being_observed"
$ws.Range("C1").Value = "type_of_quantity_observed"
$ws.Range("B1").Value = "devicetype"
$ws.Range("A1").Value = "deviceid"

# Populate the previously-blank data row 3 with a second device entry.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Device 2"
$ws.Range("C3").Value = "Temperature "
$ws.Range("D3").Value = "Celsius"

# Move the active selection to B3 (single cell) to match the saved view state.
$ws.Range("B3").Select()
